# Slide 2, "Text Placeholder 3" shape: the 4th paragraph currently reads
# "Fejlesztés stádiuma: Fejlesztés alatt" and should become
# "Fejlesztés stádiuma: DEMO" (splitting the trailing run into a
# " " run plus a new "DEMO" run).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$para = $tr.Paragraphs(4)
$para.Text = "Fejlesztés stádiuma: DEMO"
